$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows above the current row 331, shifting the
# existing rows 331..351 down to 333..353 (matches the target dimension
# A1:R353).
$ws.Range("A331:A332").EntireRow.Insert()

# New row 331: Coliflor, Primera, Feria Lagunitas de Puerto Montt
$ws.Range("A331").Value = 4
$ws.Range("B331").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C331").Value = "Los Lagos"
$ws.Range("D331").Value = 44746
$ws.Range("E331").Value = 10
$ws.Range("F331").Value = 100112008
$ws.Range("G331").Value = "Coliflor"
$ws.Range("H331").Value = "Sin especificar"
$ws.Range("I331").Value = "Primera"
$ws.Range("J331").Value = 500
$ws.Range("K331").Value = 1800
$ws.Range("L331").Value = 1800
$ws.Range("M331").Value = 1800
$ws.Range("N331").Value = '$/unidad'
$ws.Range("O331").Value = "Región del Maule"
$ws.Range("P331").Value = 1800
$ws.Range("Q331").Value = 1
$ws.Range("R331").Value = "Hortaliza"

# New row 332: Coliflor, Segunda, Feria Lagunitas de Puerto Montt
$ws.Range("A332").Value = 4
$ws.Range("B332").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C332").Value = "Los Lagos"
$ws.Range("D332").Value = 44746
$ws.Range("E332").Value = 10
$ws.Range("F332").Value = 100112008
$ws.Range("G332").Value = "Coliflor"
$ws.Range("H332").Value = "Sin especificar"
$ws.Range("I332").Value = "Segunda"
$ws.Range("J332").Value = 250
$ws.Range("K332").Value = 1500
$ws.Range("L332").Value = 1500
$ws.Range("M332").Value = 1500
$ws.Range("N332").Value = '$/unidad'
$ws.Range("O332").Value = "Región del Maule"
$ws.Range("P332").Value = 1500
$ws.Range("Q332").Value = 1
$ws.Range("R332").Value = "Hortaliza"
